$d = $word.ActiveDocument

# --- Step 1: remove the "License Information" Heading2 paragraph entirely ---
$targetIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
    if ($t -eq "License Information") {
        $targetIdx = $i
        break
    }
}
if ($targetIdx -eq 0) {
    throw "License Information paragraph not found"
}
$d.Paragraphs($targetIdx).Range.Delete()

# --- Step 2: rewrite the license-text paragraph (now immediately after the
#     "Resource: ..." Heading2 paragraph) ---
$targetIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.Contains("is based on") -and $t.Contains("unfoldingWord")) {
        $targetIdx = $i
        break
    }
}
if ($targetIdx -eq 0) {
    throw "License text paragraph not found"
}

$p = $d.Paragraphs($targetIdx)
# Clear all text in the paragraph but keep the paragraph itself (leave the
# final character in place so Word doesn't merge the paragraph away, then
# trim that leftover empty run's text too).
$clearRng = $d.Range($p.Range.Start, $p.Range.End - 1)
$clearRng.Delete()

$p = $d.Paragraphs($targetIdx)
$pos = $p.Range.Start

$text1 = "unfoldingWord® Translation Questions"
$ip1 = $d.Range($pos, $pos)
$ip1.InsertAfter($text1)
$end1 = $pos + $text1.Length
$r1 = $d.Range($pos, $end1)
$r1.Font.Bold = 1
$pos = $end1

$text2 = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. "
$ip2 = $d.Range($pos, $pos)
$ip2.InsertAfter($text2)
$end2 = $pos + $text2.Length
$r2 = $d.Range($pos, $end2)
$r2.Font.Bold = 0
$pos = $end2

$text3 = "unfoldingWord® Translation Questions"
$ip3 = $d.Range($pos, $pos)
$ip3.InsertAfter($text3)
$end3 = $pos + $text3.Length
$r3 = $d.Range($pos, $end3)
$r3.Font.Bold = 0
$pos = $end3

$text4 = " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from "
$ip4 = $d.Range($pos, $pos)
$ip4.InsertAfter($text4)
$end4 = $pos + $text4.Length
$r4 = $d.Range($pos, $end4)
$r4.Font.Bold = 0
$pos = $end4

$text5 = "unfoldingWord® Translation Questions"
$ip5 = $d.Range($pos, $pos)
$ip5.InsertAfter($text5)
$end5 = $pos + $text5.Length
$r5 = $d.Range($pos, $end5)
$r5.Font.Bold = 0
$pos = $end5

$text6 = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual"
$ip6 = $d.Range($pos, $pos)
$ip6.InsertAfter($text6)
$end6 = $pos + $text6.Length
$r6 = $d.Range($pos, $end6)
$r6.Font.Bold = 0
$pos = $end6

Write-Output "Rewritten paragraph: [$($d.Paragraphs($targetIdx).Range.Text)]"

# --- Step 3: remove the "This PDF version is provided under the same
#     license." paragraph entirely (now right after the rewritten paragraph) ---
$targetIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
    if ($t -eq "This PDF version is provided under the same license.") {
        $targetIdx = $i
        break
    }
}
if ($targetIdx -eq 0) {
    throw "'This PDF version...' paragraph not found"
}
$d.Paragraphs($targetIdx).Range.Delete()

Write-Output "Done. Paragraph count = $($d.Paragraphs.Count)"
